$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d2 = $ws.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "27.751.67"
$d2.Style = "Normal"
$ws.Range("E2").Value = "  +5.87%  "
$d3 = $ws.Range("D3")
$d3.NumberFormat = "@"
$d3.Value = "1.730.68"
$d3.Style = "Normal"
$ws.Range("E3").Value = "  +4.32%  "
$d4 = $ws.Range("D4")
$d4.NumberFormat = "@"
$d4.Value = "1.002"
$d4.Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$d5 = $ws.Range("D5")
$d5.NumberFormat = "@"
$d5.Value = "227.16"
$d5.Style = "Normal"
$ws.Range("E5").Value = "  +3.48%  "
$d6 = $ws.Range("D6")
$d6.NumberFormat = "@"
$d6.Value = "0.5435"
$d6.Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "
$d7 = $ws.Range("D7")
$d7.NumberFormat = "@"
$d7.Value = "1.002"
$d7.Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$d8 = $ws.Range("D8")
$d8.NumberFormat = "@"
$d8.Value = "0.2727"
$d8.Style = "Normal"
$ws.Range("E8").Value = "  +1.42%  "
$d9 = $ws.Range("D9")
$d9.NumberFormat = "@"
$d9.Value = "0.06666"
$d9.Style = "Normal"
$ws.Range("E9").Value = "  +4.36%  "
$d10 = $ws.Range("D10")
$d10.NumberFormat = "@"
$d10.Value = "21.83"
$d10.Style = "Normal"
$ws.Range("E10").Value = "  +5.63%  "
$d11 = $ws.Range("D11")
$d11.NumberFormat = "@"
$d11.Value = "0.07771"
$d11.Style = "Normal"
$ws.Range("E11").Value = "  +1.06%  "
$d12 = $ws.Range("D12")
$d12.NumberFormat = "@"
$d12.Value = "4.686"
$d12.Style = "Normal"
$ws.Range("E12").Value = "  +1.33%  "
$d13 = $ws.Range("D13")
$d13.NumberFormat = "@"
$d13.Value = "1.730.45"
$d13.Style = "Normal"
$ws.Range("E13").Value = "  +4.83%  "
$d14 = $ws.Range("D14")
$d14.NumberFormat = "@"
$d14.Value = "1.967.92"
$d14.Style = "Normal"
$ws.Range("E14").Value = "  +4.24%  "
$d15 = $ws.Range("D15")
$d15.NumberFormat = "@"
$d15.Value = "0.5944"
$d15.Style = "Normal"
$ws.Range("E15").Value = "  +5.26%  "
$d16 = $ws.Range("D16")
$d16.NumberFormat = "@"
$d16.Value = "0.0₅8380"
$d16.Style = "Normal"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("E17").Value = "  +4.62%  "
$d18 = $ws.Range("D18")
$d18.NumberFormat = "@"
$d18.Value = "27.740.61"
$d18.Style = "Normal"
$ws.Range("E18").Value = "  +5.93%  "
$d19 = $ws.Range("D19")
$d19.NumberFormat = "@"
$d19.Value = "225.50"
$d19.Style = "Normal"
$ws.Range("E19").Value = "  +17.49%  "
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +4.12%  "
$d23 = $ws.Range("D23")
$d23.NumberFormat = "@"
$d23.Value = "6.198"
$d23.Style = "Normal"
$ws.Range("E23").Value = "  +3.23%  "
$d24 = $ws.Range("D24")
$d24.NumberFormat = "@"
$d24.Value = "1.003"
$d24.Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$d25 = $ws.Range("D25")
$d25.NumberFormat = "@"
$d25.Value = "147.38"
$d25.Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$d26 = $ws.Range("D26")
$d26.NumberFormat = "@"
$d26.Value = "1.724"
$d26.Style = "Normal"
$ws.Range("E26").Value = "  +12.93%  "
$d27 = $ws.Range("D27")
$d27.NumberFormat = "@"
$d27.Value = "0.1248"
$d27.Style = "Normal"
$ws.Range("E27").Value = "  +3.55%  "
$d28 = $ws.Range("D28")
$d28.NumberFormat = "@"
$d28.Value = "7.454"
$d28.Style = "Normal"
$ws.Range("E28").Value = "  +2.09%  "
$d29 = $ws.Range("D29")
$d29.NumberFormat = "@"
$d29.Value = "17.02"
$d29.Style = "Normal"
$ws.Range("E29").Value = "  +5.86%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +2.48%  "
$d32 = $ws.Range("D32")
$d32.NumberFormat = "@"
$d32.Value = "3.654"
$d32.Style = "Normal"
$ws.Range("E32").Value = "  +4.48%  "
$d33 = $ws.Range("D33")
$d33.NumberFormat = "@"
$d33.Value = "3.506"
$d33.Style = "Normal"
$ws.Range("E33").Value = "  +2.99%  "
$d34 = $ws.Range("D34")
$d34.NumberFormat = "@"
$d34.Value = "1.673"
$d34.Style = "Normal"
$ws.Range("E34").Value = "  +5.68%  "
$d35 = $ws.Range("D35")
$d35.NumberFormat = "@"
$d35.Value = "0.9716"
$d35.Style = "Normal"
$ws.Range("E35").Value = "  +1.90%  "
$d36 = $ws.Range("D36")
$d36.NumberFormat = "@"
$d36.Value = "2.846"
$d36.Style = "Normal"
$ws.Range("E36").Value = "  +1.91%  "
$d37 = $ws.Range("D37")
$d37.NumberFormat = "@"
$d37.Value = "2.435"
$d37.Style = "Normal"
$ws.Range("E37").Value = "  +1.29%  "
$d38 = $ws.Range("D38")
$d38.NumberFormat = "@"
$d38.Value = "0.5977"
$d38.Style = "Normal"
$ws.Range("E38").Value = "  +3.42%  "
$d39 = $ws.Range("D39")
$d39.NumberFormat = "@"
$d39.Value = "0.01670"
$d39.Style = "Normal"
$ws.Range("E39").Value = "  +3.94%  "
$d40 = $ws.Range("D40")
$d40.NumberFormat = "@"
$d40.Value = "5.905"
$d40.Style = "Normal"
$ws.Range("E40").Value = "  -1.36%  "
$d41 = $ws.Range("D41")
$d41.NumberFormat = "@"
$d41.Value = "0.8626"
$d41.Style = "Normal"
$ws.Range("E41").Value = "  +3.16%  "
$d42 = $ws.Range("D42")
$d42.NumberFormat = "@"
$d42.Value = "1.048.69"
$d42.Style = "Normal"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("E43").Value = "  -0.13%  "
$d44 = $ws.Range("D44")
$d44.NumberFormat = "@"
$d44.Value = "101.51"
$d44.Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$d45 = $ws.Range("D45")
$d45.NumberFormat = "@"
$d45.Value = "1.872.71"
$d45.Style = "Normal"
$ws.Range("E45").Value = "  +4.13%  "
$ws.Range("E46").Value = "  +9.22%  "
$d47 = $ws.Range("D47")
$d47.NumberFormat = "@"
$d47.Value = "59.54"
$d47.Style = "Normal"
$ws.Range("E47").Value = "  +1.58%  "
$d48 = $ws.Range("D48")
$d48.NumberFormat = "@"
$d48.Value = "8.254"
$d48.Style = "Normal"
$ws.Range("E48").Value = "  +2.03%  "
$d49 = $ws.Range("D49")
$d49.NumberFormat = "@"
$d49.Value = "0.4430"
$d49.Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "
$d50 = $ws.Range("D50")
$d50.NumberFormat = "@"
$d50.Value = "0.05328"
$d50.Style = "Normal"
$ws.Range("E50").Value = "  -0.32%  "
$d51 = $ws.Range("D51")
$d51.NumberFormat = "@"
$d51.Value = "0.9992"
$d51.Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
